# Fruta / hortaliza, semanal
# Weekly refresh of the "Ajo" (garlic) sheet: 3 new price-report rows are
# inserted above the current last block of rows (old rows 364-370 shift
# down to 367-373), and the 3 new rows (364-366) are populated with this
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing rows 364:370 down by 3 rows (same as inserting 3 new
# blank rows above row 364), carrying formatting (e.g. the date style on
# column D) down with them.
$ws.Range("A364:R366").Insert()

# New data for rows 364-366 (columns A through R).
$newRows = @(
    @(364, @(3, 'Femacal de La Calera', 'Coquimbo', 44595, 5, 100112003, 'Ajo', 'Chino', '1a (cosecha)', 73, 16000, 16500, 16260, '$/caja 10 kilos', 'Llay Llay', 1626, 10, 'Hortaliza')),
    @(365, @(3, 'Femacal de La Calera', 'Coquimbo', 44595, 5, 100112003, 'Ajo', 'Chino', '1a (cosecha)', 80, 6000, 6000, 6000, '$/trenza 50 unidades', 'Llay Llay', 1200, 5, 'Hortaliza')),
    @(366, @(3, 'Femacal de La Calera', 'Coquimbo', 44595, 5, 100112003, 'Ajo', 'Chino', '2a (guarda)', 85, 4000, 4000, 4000, '$/trenza 50 unidades', 'Llay Llay', 800, 5, 'Hortaliza'))
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
